$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-modified date serial for each row.
# Bump every value in C2:C43 from 45774 to 45775 (one day later),
# matching the "Automatic update of files." commit.
$ws.Range("C2:C43").Value = 45775
